$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3049")

# Ativação: 01/01/2018 -> 01/01/2021
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

# Objetivos:
$novoObjetivo = "Abordar os princípios básicos da termodinâmica de forma que os estudantes e futuros engenheiros tenham um entendimento claro e sólido sobre estes princípios. Apresentar diversos exemplos de engenharia do mundo real e de como a termodinâmica é aplicada na prática de engenharia. Enfatizar a compreensão da termodinâmica baseada na Física e em argumentos físicos, buscando incentivar o entendimento mais profundo da termodinâmica."
$ws.Range("B10").Value = $novoObjetivo
$ws.Range("C10").Value = $novoObjetivo

# Programa resumido:
$novoProgramaResumido = "1. Termodinâmica e Energia. 2. Importância das unidades e análise dimensional.3. Sistemas e volumes de controle. 4. Equipamentos domésticos e a Termodinâmica. 5. Propriedades de um sistema: estados termodinâmicos e equilíbrio. 6. Eficiência na conversão de energia. 7. Processos e ciclos térmicos. 8. Termodinâmica e o meio ambiente."
$ws.Range("B14").Value = $novoProgramaResumido
$ws.Range("C14").Value = $novoProgramaResumido

# Programa:
$novoPrograma = "1. Termodinâmica e Energia: formas de energia e transferência de energia por calor e trabalho; formas mecânicas de trabalho. 2. Sistema de Unidades e Análise Dimensional: importância na engenharia de máquinas. 3. Sistemas e volumes de controle: dispositivos ativos e passivos. 4. Propriedades de um sistema. Estados e equilíbrio: diagramas de propriedades para processos com mudança de fase; equilíbrio de estado do gás ideal; fator de compressibilidade; pressão de vapor e pressão de equilíbrio; calores específicos. 5. Balanço de energia em sistemas fechados e em volumes de controle: trabalho de fluxo e energia de escoamento de um fluido; regime permanente e transiente. 6. Máquinas térmicas e refrigeradores e a 2ª. Lei da Termodinâmica: princípios e ciclos de Carnot; entropia e variação de entropia em sólidos, líquidos e gases. 7. Eficiência na conversão de energia. Eficiência térmica. Eficiência de máquinas. Eficiência isoentrópica em dispositivos com escoamento em regime permanente. Balanço de entropia. 8. Processo e ciclos: Ciclos de potência a gás: Otto, Diesel, Stirling, Ericsson, Brayton e suas variações. Ciclos de potência a vapor e ciclos combinados gás-vapor: Rankine ideal; afastamento da condição ideal; eficiência do ciclo Rankine com e sem modificações; cogeração. Ciclos de refrigeração e sistemas de bombas de calor: sistemas a gás e por absorção. 9. Economia de energia: benefícios ao meio ambiente."
$ws.Range("B16").Value = $novoPrograma
$ws.Range("C16").Value = $novoPrograma
